$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J - copy formatting from H1 (bold/border/centered header style)
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$data = @(
    @(2,8,8),
    @(3,9,9),
    @(4,9,9),
    @(5,7,7),
    @(6,8,8),
    @(7,9,9),
    @(8,9,9),
    @(9,9,9),
    @(10,9,9),
    @(11,9,9),
    @(12,9,9),
    @(13,8,8),
    @(14,9,9),
    @(15,9,9),
    @(16,9,9),
    @(17,9,9),
    @(18,9,9),
    @(19,8,8),
    @(20,9,9),
    @(21,8,9),
    @(22,9,9),
    @(23,10,10),
    @(24,9,9),
    @(25,9,9),
    @(26,9,9),
    @(27,9,9),
    @(28,8,8),
    @(29,8,8),
    @(30,8,8),
    @(31,9,9),
    @(32,9,9),
    @(33,8,8),
    @(34,8,8),
    @(35,7,8),
    @(36,8,8),
    @(37,8,8),
    @(38,8,8),
    @(39,8,8),
    @(40,7,7),
    @(41,7,8),
    @(42,9,9),
    @(43,8,8),
    @(44,7,7),
    @(45,8,8),
    @(46,8,8),
    @(47,7,7),
    @(48,7,8),
    @(49,8,8),
    @(50,8,8),
    @(51,7,7),
    @(52,8,8),
    @(53,9,9),
    @(54,5,5),
    @(55,5,5),
    @(56,6,6)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($row, 9).Value = $iVal
    $ws.Cells.Item($row, 10).Value = $jVal
}
